$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row values (rows 2-25) for columns B..N (index 0-42: B=0 .. N=12)
$data = @(
    @(0.8258382392119188, 0.169130870859135, 0, 0.4201928946818256, 2.188402845309867, 0.6932782575092489, 0.7473974425159327, 0, 0.04093914717267033, 0, 0, 0.5224473424993192, 1.190615592527177),
    @(0.7345745284280838, 0.1475884978857778, 0, 0.4182430313368641, 2.166182142814748, 0.6789583914570159, 0.7463789717493938, 0, 0.04159194831786817, 0, 0, 0.4941220946971541, 1.208097573569571),
    @(0.6786953740740955, 0.1343205559346927, 0, 0.4171838556245469, 2.153872985124437, 0.6707938444422865, 0.7462467883566006, 0, 0.04201820398814249, 0, 0, 0.4769575867037616, 1.219378047094391),
    @(0.6559636986867474, 0.1289034637228212, 0, 0.4167870263868743, 2.149191853536038, 0.667623692121694, 0.7463166238134988, 0, 0.04219828478950305, 0, 0, 0.4700201826739629, 1.22411217004081),
    @(0.6521915100061335, 0.1280033369282307, 0, 0.4167232370038079, 2.148434767796033, 0.6671067441366745, 0.7463356822208027, 0, 0.0422285717928057, 0, 0, 0.4688716945816722, 1.224906551372982),
    @(0.678388646642901, 0.1342475408653456, 0, 0.4171783628624368, 2.153808498354806, 0.6707504563339768, 0.7462472297332141, 0, 0.04202060681924458, 0, 0, 0.47686379440011, 1.219441337695391),
    @(0.7943378803104224, 0.1617115121189272, 0, 0.4194919764214653, 2.180463803662747, 0.6882098809967232, 0.7469437294120382, 0, 0.04115894404298714, 0, 0, 0.5126336765949091, 1.196529856138355),
    @(1.022972125688511, 0.2152482490811849, 0, 0.4251213910440015, 2.243356318210473, 0.7274745104736553, 0.7522377330166847, 0, 0.03967178027370721, 0, 0, 0.584581151045434, 1.155943701018495),
    @(1.191750234745257, 0.2543960113758601, 0, 0.4299202405926295, 2.296093396280284, 0.7594524398460578, 0.7585447336814184, 0, 0.0387036549152171, 0, 0, 0.638545313450777, 1.128781812921993),
    @(1.268714011767372, 0.2721673122339325, 0, 0.4322467291632321, 2.321515224393977, 0.7746941806902328, 0.761943954729901, 0, 0.03829048196739926, 0, 0, 0.6633363860717623, 1.117004060657731),
    @(1.297885210333106, 0.2788915925816582, 0, 0.4331482760470848, 2.331348505092819, 0.7805668407144992, 0.7633077686868432, 0, 0.03813795957469956, 0, 0, 0.6727589903848923, 1.112627506889167),
    @(1.291601484210901, 0.2774436349208429, 0, 0.4329531987302531, 2.329221532260277, 0.7792975523681207, 0.7630106340955365, 0, 0.03817063254949637, 0, 0, 0.6707281195456574, 1.113566361433071),
    @(1.271113409110342, 0.2727206294032669, 0, 0.4323204884151366, 2.322320069737202, 0.7751752989927638, 0.7620546191623134, 0, 0.03827785482718937, 0, 0, 0.6641108933039703, 1.116642325449924),
    @(1.258567346278312, 0.2698269593976761, 0, 0.4319356099023963, 2.318119650466144, 0.7726634762342712, 0.7614790193810848, 0, 0.0383440449482304, 0, 0, 0.6600621751863542, 1.118537314576177),
    @(1.186724206544, 0.2532338683868716, 0, 0.4297710796602843, 2.29446088811882, 0.7584704256681505, 0.7583332858118013, 0, 0.03873120653355677, 0, 0, 0.6369300289875213, 1.12956317646392),
    @(1.142698324996786, 0.2430450308840761, 0, 0.4284799028815058, 2.28031416271439, 0.7499421397851336, 0.7565395226142186, 0, 0.03897570750487134, 0, 0, 0.6228012392765407, 1.136475493966051),
    @(1.117393288781841, 0.2371811767381757, 0, 0.4277507604217732, 2.272312031974238, 0.7451022332793968, 0.7555576739601975, 0, 0.03911889992819795, 0, 0, 0.6146975650725039, 1.14050569496065),
    @(1.10882844260658, 0.2351951738763489, 0, 0.4275062076489959, 2.269625758886946, 0.7434747146705263, 0.7552337919853187, 0, 0.03916782179126521, 0, 0, 0.6119577250065618, 1.14187959033514),
    @(1.147383145687058, 0.244130012528899, 0, 0.4286159533409375, 2.281806160394396, 0.7508432207164617, 0.7567253070868958, 0, 0.03894941465602564, 0, 0, 0.6243029108405267, 1.135734031966424),
    @(1.277130529304827, 0.2741080349309755, 0, 0.4325057734724709, 2.324341583923825, 0.7763833569271696, 0.76233334198119, 0, 0.03824625401938953, 0, 0, 0.6660535904415639, 1.115736573990031),
    @(1.362083489119527, 0.2936693774267383, 0, 0.4351677788273491, 2.353345434801895, 0.7936641698746598, 0.7664451902732878, 0, 0.03780965715487916, 0, 0, 0.6935426514670979, 1.103153472499157),
    @(1.316728269627959, 0.2832319508053729, 0, 0.4337360806810935, 2.337755074950053, 0.7843868529122631, 0.7642096257334003, 0, 0.03804056925253185, 0, 0, 0.6788527212869013, 1.109824717198466),
    @(1.145265119472526, 0.2436395115062453, 0, 0.4285544038934574, 2.281131220217674, 0.750435645499067, 0.7566411600465983, 0, 0.03896129347596933, 0, 0, 0.6236239454537014, 1.136069071841025),
    @(0.9609823812872378, 0.2007984392826643, 0, 0.4234818519158807, 2.225199485342813, 0.7163074592028522, 0.750382684294479, 0, 0.04005230891283684, 0, 0, 0.5649239122673677, 1.166457452264634)
)

$rowCount = 24
$colCount = 13
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$ws.Range("B2:N25").Value = $arr
